# Edit script: update ultrasound machine brand/model mentioned in the report
# template and tidy up a spell-check artifact ("ecogenicidad" run split).

$d = $word.ActiveDocument

# 1) Replace the ultrasound equipment brand/model text.
#    "MEDISONIC MODELO H60 " -> "MINDRAY MODELO DC – N3 "
$d.Content.Find.Execute(
    "MEDISONIC MODELO H60 ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "MINDRAY MODELO DC " + [char]0x2013 + " N3 ",
    2
)

# 2) Merge the runs around "ecogenicidad" that were split apart (spell-check
#    artifact): re-typing the same text via Find/Replace collapses the
#    separate "ecogenicidad" run (wrapped in proofErr) and the following
#    lone-space run into the preceding run.
$d.Content.Find.Execute(
    "paredes lisas y la ecogenicidad ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "paredes lisas y la ecogenicidad ",
    2
)

# 3) Tune the run formatting on the newly-inserted brand/model text: drop the
#    w:cs="Tahoma" override on rFonts and bump szCs from 18 to 20.
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "MINDRAY MODELO DC " + [char]0x2013 + " N3 "
$rng.Find.Forward = $true
$rng.Find.Wrap = 1
$rng.Find.Execute()
if ($rng.Find.Found) {
    $rng.Font.NameOther = $rng.Font.NameAscii
    $rng.Font.SizeBi = 10
}
